$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values per row (only cells that actually change are listed;
# unchanged cells are left untouched).
$updates = @{
    2  = @{ G=20.187688; H=60.563064; I=0.4610390972738712; J=0.4610390972738712;
            M=10.00528166666667; N=30.015845; O=0.2411298117910886; P=0.2411298117910886;
            Q=201.9835046387867; R=1817.85154174908; S=0.1111702707539819; T=0.1111702707539819 }
    3  = @{ G=20.187688; H=60.563064; I=0.4610390972738712; J=0.4610390972738712;
            O=0.6095920530221616; P=0.6095920530221617;
            Q=510.6276090657973; R=4595.648481592176; S=0.2810457698306632; T=0.2810457698306633 }
    4  = @{ G=20.187688; H=60.563064; I=0.4610390972738712; J=0.4610390972738712;
            M=6.194048666666666;
            O=0.1492781351867498; P=0.1492781351867498;
            Q=125.0435219394826; R=1125.391697455344; S=0.06882305668922604; T=0.06882305668922604 }
    5  = @{ I=0.4155251476643568; J=0.4155251476643568;
            M=10.00528166666667; N=30.015845; O=0.2411298117910886; P=0.2411298117910886;
            Q=182.0436186151467; R=1638.39256753632; S=0.1001955006507706; T=0.1001955006507706 }
    6  = @{ I=0.4155251476643568; J=0.4155251476643568;
            O=0.6095920530221616; P=0.6095920530221617;
            Q=460.2182632951894; R=4141.964369656705; S=0.2533008278470521; T=0.2533008278470522 }
    7  = @{ I=0.4155251476643568; J=0.4155251476643568;
            M=6.194048666666666;
            O=0.1492781351867498; P=0.1492781351867498;
            S=0.06202881916653404; T=0.06202881916653403 }
    8  = @{ I=0.123435755061772; J=0.123435755061772;
            M=10.00528166666667; N=30.015845; O=0.2411298117910886; P=0.2411298117910886;
            Q=54.07781368767778; R=486.7003231891; S=0.029764040386336; T=0.02976404038633599 }
    9  = @{ I=0.123435755061772; J=0.123435755061772;
            O=0.6095920530221616; P=0.6095920530221617;
            S=0.07524545534444628; T=0.07524545534444629 }
    10 = @{ I=0.123435755061772; J=0.123435755061772;
            M=6.194048666666666;
            O=0.1492781351867498; P=0.1492781351867498;
            S=0.01842625933098974; T=0.01842625933098974 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
